{"js": "// Update the date line and the 25 division problems in the practice table.\n// The table has 20 rows total; only rows 0, 4, 8, 12 and 16 (0-based) hold\n// the 5 problems per row, the rest are blank answer rows.\n\nconst body = context.document.body;\n\n// 1) Update the date/title paragraph.\nconst titlePara = body.paragraphs.getFirst();\ntitlePara.insertText(\"2026-01-13 Tuesday\", Word.InsertLocation.replace);\n\n// 2) Update the division problems, preserving cell run formatting by\n//    setting Table.getCell(row, col).value.\nconst table = body.tables.getFirst();\n\nconst newValues = {\n  0: [\"430\u00f78=\", \"904\u00f78=\", \"527\u00f72=\", \"794\u00f78=\", \"383\u00f78=\"],\n  4: [\"846\u00f72=\", \"717\u00f72=\", \"598\u00f79=\", \"514\u00f78=\", \"476\u00f72=\"],\n  8: [\"564\u00f72=\", \"625\u00f73=\", \"762\u00f77=\", \"390\u00f78=\", \"696\u00f74=\"],\n  12: [\"814\u00f75=\", \"333\u00f77=\", \"192\u00f72=\", \"204\u00f74=\", \"856\u00f77=\"],\n  16: [\"542\u00f74=\", \"813\u00f73=\", \"187\u00f72=\", \"478\u00f79=\", \"449\u00f73=\"],\n};\n\nfor (const rowIndex of Object.keys(newValues)) {\n  const row = parseInt(rowIndex, 10);\n  const values = newValues[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice table.\n# The table has 20 rows total; only rows 1, 5, 9, 13 and 17 (1-based, as\n# used by Word COM's Cell(row, col)) hold the 5 problems per row - the rest\n# are blank answer rows.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/title paragraph.\n$d.Paragraphs.Item(1).Range.Text = \"2026-01-13 Tuesday\"\n\n# 2) Update the division problems, preserving cell run formatting by\n#    setting the text of each cell's Range.\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1 = @(\"430\u00f78=\", \"904\u00f78=\", \"527\u00f72=\", \"794\u00f78=\", \"383\u00f78=\")\n    5 = @(\"846\u00f72=\", \"717\u00f72=\", \"598\u00f79=\", \"514\u00f78=\", \"476\u00f72=\")\n    9 = @(\"564\u00f72=\", \"625\u00f73=\", \"762\u00f77=\", \"390\u00f78=\", \"696\u00f74=\")\n    13 = @(\"814\u00f75=\", \"333\u00f77=\", \"192\u00f72=\", \"204\u00f74=\", \"856\u00f77=\")\n    17 = @(\"542\u00f74=\", \"813\u00f73=\", \"187\u00f72=\", \"478\u00f79=\", \"449\u00f73=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
